# Add header row to the worksheet: Company Name, Customer ID, Invoice Number, Total amount
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Company Name"
$ws.Range("B1").Value = "Customer ID"
$ws.Range("C1").Value = "Invoice Number"
$ws.Range("D1").Value = "Total amount"

# Give the new header cells an explicit (black) font color, same as the
# shipped workbook's header styling.
$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Color = 0
